$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert a new "Meta description" paragraph right after the first
# (Heading1 title) paragraph.
#   <w:p>
#     <w:r/>
#     <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#     <w:r><w:t>: Find out about the Dark Mystic slot game by Felix Gaming.
#               Learn about its features, pros and cons, and play for free
#               in 2021.</w:t></w:r>
#   </w:p>
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

# Put the new paragraph on the Normal style *before* anything else touches it
# (doing this first keeps the paragraph's <w:pPr> absent, matching the rest
# of the body-text paragraphs in the document).
$d.Paragraphs.Item(2).Style = "Normal"

# Borrow the run layout (empty run + bold run) of the existing bold title
# paragraph near the end of the document ("Play Dark Mystic Slot for
# Free..."), which already has the <w:r/><w:r><w:rPr><w:b/></w:rPr>...
# shape we need for the new paragraph.
$boldTemplatePara = $d.Paragraphs.Item(51)
$d.Paragraphs.Item(2).Range.FormattedText = $boldTemplatePara.Range.FormattedText

# Swap the copied bold text for "Meta description".
$d.Paragraphs.Item(2).Range.Find.Execute(
    "Play Dark Mystic Slot for Free - Review & Pros/Cons - 2021",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Meta description", 2) | Out-Null

# Append the (non-bold) rest of the sentence right after the bold run, before
# the paragraph mark.
$metaPara = $d.Paragraphs.Item(2)
$endPos = $metaPara.Range.End
$appendPoint = $d.Range($endPos - 1, $endPos - 1)
$appendPoint.InsertAfter(": Find out about the Dark Mystic slot game by Felix Gaming. Learn about its features, pros and cons, and play for free in 2021.") | Out-Null

# ---------------------------------------------------------------------------
# Change 2: near the end of the document, drop the bold
# "Play Dark Mystic Slot for Free..." paragraph entirely, and replace the
# text of the remaining (formerly italic "Find out about...") paragraph with
# the DALLE prompt, keeping its existing <w:r/><w:r><w:rPr><w:i/>.../w:r>
# shape intact.
# ---------------------------------------------------------------------------

$boldTailPara = $d.Paragraphs.Item(51)
$boldTailPara.Range.Delete() | Out-Null

$dallePrompt = "Prompt for DALLE: Create a cartoon-style image featuring a happy Maya warrior with glasses for the game " + [char]34 + "Dark Mystic" + [char]34 + ". The image should be visually appealing and capture the essence of the game" + [char]0x2019 + "s fantasy genre. Use vibrant colors and include fiery elements like a dragon or cascading wins to add excitement. The character should be happy and enthusiastic to reflect the lively gameplay of the slot game."

$italicTailPara = $d.Paragraphs.Item(51)
$italicTailPara.Range.Find.Execute(
    "Find out about the Dark Mystic slot game by Felix Gaming. Learn about its features, pros and cons, and play for free in 2021.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $dallePrompt, 2) | Out-Null
